$d = $word.ActiveDocument

# --- Split 1: the run containing "{m" must become two runs: "{" and "m",
#     each keeping the original run's formatting (w:rPr with w:lang en-US).
#
# Locate the field opening "{m:'" (there's also an unrelated "{m:commentblock}"
# literal text run earlier in the document, so search from where the field
# construct starts, i.e. after the introductory paragraphs).
$search1 = $d.Content.Duplicate
$found1 = $search1.Find.Execute("{m:'", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate the '{m:''' field opening text"
}

$openBraceStart = $search1.Start
$openBraceEnd = $openBraceStart + 1

# Toggling a character formatting property on the first character ("{") and
# back to its original value forces Word to split the run in two while
# leaving the effective formatting (and therefore the serialized w:rPr)
# unchanged on both halves.
$rOpenBrace = $d.Range($openBraceStart, $openBraceEnd)
$rOpenBrace.Font.Bold = 1
$rOpenBrace.Font.Bold = 0

# --- Split 2: the run containing ")}" must become two runs: ")" and "}",
#     where the new "}" run has NO run properties at all (no w:rPr).
$search2 = $d.Content.Duplicate
$found2 = $search2.Find.Execute(")}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the ')}' field closing text"
}

$closeParenStart = $search2.Start
$closeParenEnd = $closeParenStart + 1
$closeBraceStart = $closeParenEnd
$closeBraceEnd = $closeBraceStart + 1

# Remove the trailing "}" then re-insert it with InsertAfter: text inserted
# this way lands in a brand new run that carries no formatting (no w:rPr),
# exactly like a plain, unformatted run typed at the end of the story.
$rCloseBrace = $d.Range($closeBraceStart, $closeBraceEnd)
$rCloseBrace.Text = ""

$rCloseParen = $d.Range($closeParenStart, $closeParenEnd)
$rCloseParen.InsertAfter("}")
